$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = 44490
$ws.Range("K6").Value = 450
$ws.Range("L6").Value = 480
$ws.Range("M6").Value = 465
$ws.Range("P6").Value = 465

# Row 7
$ws.Range("D7").Value = 44322
$ws.Range("I7").Value = "1a (cosecha)"
$ws.Range("K7").Value = 350
$ws.Range("L7").Value = 400
$ws.Range("M7").Value = 375
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 375

# Row 8
$ws.Range("D8").Value = 44665
$ws.Range("K8").Value = 400
$ws.Range("L8").Value = 420
$ws.Range("M8").Value = 410
$ws.Range("O8").Value = "Región de O'Higgins"
$ws.Range("P8").Value = 410
